$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3605800351893436
$ws.Range("C2").Value = 0.03088557933001823
$ws.Range("D2").Value = 0.1582768291396235
$ws.Range("E2").Value = 0.1566086668042352
$ws.Range("F2").Value = 1.732319223695484
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1935279873024527
$ws.Range("K2").Value = 0.3142421109734244
$ws.Range("M2").Value = 0.2122778647614751
$ws.Range("O2").Value = 4.396756074374139

$ws.Range("B3").Value = 0.3283216481623867
$ws.Range("C3").Value = 0.02738043922929023
$ws.Range("D3").Value = 0.1550313151745257
$ws.Range("E3").Value = 0.1554378733549875
$ws.Range("F3").Value = 1.737500715248814
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1934081094813109
$ws.Range("K3").Value = 0.2810876666782747
$ws.Range("M3").Value = 0.2020710694219332
$ws.Range("O3").Value = 4.41983796681879

$ws.Range("B4").Value = 0.3085938214295254
$ws.Range("C4").Value = 0.02521660268811843
$ws.Range("D4").Value = 0.1530991578740384
$ws.Range("E4").Value = 0.154786353056064
$ws.Range("F4").Value = 1.741507268733137
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1934220207013482
$ws.Range("K4").Value = 0.2607568303609611
$ws.Range("M4").Value = 0.195887129392144
$ws.Range("O4").Value = 4.43613697541997

$ws.Range("B5").Value = 0.3005748792305099
$ws.Range("C5").Value = 0.02433193968214198
$ws.Range("D5").Value = 0.1523271114381117
$ws.Range("E5").Value = 0.1545378302330285
$ws.Range("F5").Value = 1.743347577505176
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.193449718883457
$ws.Range("K5").Value = 0.2524788782477003
$ws.Range("M5").Value = 0.1933881755924034
$ws.Range("O5").Value = 4.443313832954047

$ws.Range("B6").Value = 0.2992445799039274
$ws.Range("C6").Value = 0.02418486952095122
$ws.Range("D6").Value = 0.1521998417023838
$ws.Range("E6").Value = 0.1544975897796128
$ws.Range("F6").Value = 1.743665702032821
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1934556493204056
$ws.Range("K6").Value = 0.2511047686835752
$ws.Range("M6").Value = 0.1929745024975134
$ws.Range("O6").Value = 4.444537853624993

$ws.Range("B7").Value = 0.3084855923437999
$ws.Range("C7").Value = 0.02520468341437265
$ws.Range("D7").Value = 0.1530886836295764
$ws.Range("E7").Value = 0.1547829325965857
$ws.Range("F7").Value = 1.741531247024788
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1934223050232546
$ws.Range("K7").Value = 0.260645161934363
$ws.Range("M7").Value = 0.1958533421725619
$ws.Range("O7").Value = 4.436231599172316

$ws.Range("B8").Value = 0.3494412181356097
$ws.Range("C8").Value = 0.02967945817366058
$ws.Range("D8").Value = 0.1571452475170787
$ws.Range("E8").Value = 0.1561910222056433
$ws.Range("F8").Value = 1.733934627502158
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.193468502326283
$ws.Range("K8").Value = 0.3028053005196796
$ws.Range("M8").Value = 0.2087414135881502
$ws.Range("O8").Value = 4.404273395979743

$ws.Range("B9").Value = 0.4303657638029676
$ws.Range("C9").Value = 0.03836005946547516
$ws.Range("D9").Value = 0.1655778204644776
$ws.Range("E9").Value = 0.1594851904314964
$ws.Range("F9").Value = 1.72558015008282
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1942527371958462
$ws.Range("K9").Value = 0.3856725126036906
$ws.Range("M9").Value = 0.2346683482246092
$ws.Range("O9").Value = 4.35847440088844

$ws.Range("B10").Value = 0.49017827233709
$ws.Range("C10").Value = 0.04467825869865294
$ws.Range("D10").Value = 0.1720609437847571
$ws.Range("E10").Value = 0.1622287232872139
$ws.Range("F10").Value = 1.723426786165831
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1952510818985331
$ws.Range("K10").Value = 0.4466563364004799
$ws.Range("M10").Value = 0.2541097358512587
$ws.Range("O10").Value = 4.335110384855483

$ws.Range("B11").Value = 0.5174633051150295
$ws.Range("C11").Value = 0.04753930817706475
$ws.Range("D11").Value = 0.1750720700491257
$ws.Range("E11").Value = 0.1635467263090931
$ws.Range("F11").Value = 1.723311738136687
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1957967994368843
$ws.Range("K11").Value = 0.4744185936270071
$ws.Range("M11").Value = 0.2630384031907553
$ws.Range("O11").Value = 4.326714858245595

$ws.Range("B12").Value = 0.5278059967516811
$ws.Range("C12").Value = 0.04862078185092855
$ws.Range("D12").Value = 0.1762211348655001
$ws.Range("E12").Value = 0.1640558447367546
$ws.Range("F12").Value = 1.723392402451523
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1960165987542482
$ws.Range("K12").Value = 0.4849339962837007
$ws.Range("M12").Value = 0.2664315005847158
$ws.Range("O12").Value = 4.323856747531778

$ws.Range("B13").Value = 0.5255780564370127
$ws.Range("C13").Value = 0.04838795446892163
$ws.Range("D13").Value = 0.1759732724673597
$ws.Range("E13").Value = 0.1639457520059757
$ws.Range("F13").Value = 1.723369506162172
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1959686765854372
$ws.Range("K13").Value = 0.4826692156169656
$ws.Range("M13").Value = 0.2657002051737649
$ws.Range("O13").Value = 4.324458012003447

$ws.Range("B14").Value = 0.518313998112518
$ws.Range("C14").Value = 0.04762832093709335
$ws.Range("D14").Value = 0.1751664281058822
$ws.Range("E14").Value = 0.1635884112180896
$ws.Range("F14").Value = 1.723315885177399
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1958146190644285
$ws.Range("K14").Value = 0.4752836561012543
$ws.Range("M14").Value = 0.2633173157468391
$ws.Range("O14").Value = 4.326473284535183

$ws.Range("B15").Value = 0.5138658979623472
$ws.Range("C15").Value = 0.04716276870189517
$ws.Range("D15").Value = 0.1746733582301943
$ws.Range("E15").Value = 0.1633708332774866
$ws.Range("F15").Value = 1.723299216647717
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1957219660662233
$ws.Range("K15").Value = 0.4707600904852995
$ws.Range("M15").Value = 0.261859285442533
$ws.Range("O15").Value = 4.327749513144852

$ws.Range("B16").Value = 0.488396619122625
$ws.Range("C16").Value = 0.04449101340864559
$ws.Range("D16").Value = 0.1718653981255329
$ws.Range("E16").Value = 0.1621439929404929
$ws.Range("F16").Value = 1.723451693864021
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1952172588766174
$ws.Range("K16").Value = 0.4448423758107367
$ws.Range("M16").Value = 0.2535279157281565
$ws.Range("O16").Value = 4.33570398419306

$ws.Range("B17").Value = 0.4727911635941666
$ws.Range("C17").Value = 0.0428485753280512
$ws.Range("D17").Value = 0.1701586026833439
$ws.Range("E17").Value = 0.1614092548985511
$ws.Range("F17").Value = 1.723766597907783
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1949310759260214
$ws.Range("K17").Value = 0.4289475709114186
$ws.Range("M17").Value = 0.24843845673999
$ws.Range("O17").Value = 4.341155694519927

$ws.Range("B18").Value = 0.4638225031433478
$ws.Range("C18").Value = 0.04190265572026419
$ws.Range("D18").Value = 0.1691827326041135
$ws.Range("E18").Value = 0.1609932400281764
$ws.Range("F18").Value = 1.724029101111853
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1947750907877079
$ws.Range("K18").Value = 0.4198072501794741
$ws.Range("M18").Value = 0.245519115555652
$ws.Range("O18").Value = 4.344501543058954

$ws.Range("B19").Value = 0.4607871177723553
$ws.Range("C19").Value = 0.0415821735764581
$ws.Range("D19").Value = 0.168853324343246
$ws.Range("E19").Value = 0.1608535170959833
$ws.Range("F19").Value = 1.724131959725071
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1947237579266954
$ws.Range("K19").Value = 0.4167128448164021
$ws.Range("M19").Value = 0.2445320527966572
$ws.Range("O19").Value = 4.345670486379646

$ws.Range("B20").Value = 0.4744516509295522
$ws.Range("C20").Value = 0.04302354364548933
$ws.Range("D20").Value = 0.1703396908594925
$ws.Range("E20").Value = 0.1614867875779602
$ws.Range("F20").Value = 1.723724654072711
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1949606486193289
$ws.Range("K20").Value = 0.4306394017467881
$ws.Range("M20").Value = 0.2489794135828944
$ws.Range("O20").Value = 4.340553598888448

$ws.Range("B21").Value = 0.5204473486343488
$ws.Range("C21").Value = 0.04785149690884793
$ws.Range("D21").Value = 0.1754031792683861
$ws.Range("E21").Value = 0.1636930992705459
$ws.Range("F21").Value = 1.72332826408045
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1958595128142946
$ws.Range("K21").Value = 0.4774529117318025
$ws.Range("M21").Value = 0.2640169036979714
$ws.Range("O21").Value = 4.325872635714461

$ws.Range("B22").Value = 0.5505688328302369
$ws.Range("C22").Value = 0.05099548880556881
$ws.Range("D22").Value = 0.1787638113196408
$ws.Range("E22").Value = 0.1651934328515168
$ws.Range("F22").Value = 1.723793269801277
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1965235998342649
$ws.Range("K22").Value = 0.5080622301991582
$ws.Range("M22").Value = 0.273914666431871
$ws.Range("O22").Value = 4.318149356316411

$ws.Range("B23").Value = 0.5344870431335096
$ws.Range("C23").Value = 0.04931853833248567
$ws.Range("D23").Value = 0.1769655091254805
$ws.Range("E23").Value = 0.1643873480144968
$ws.Range("F23").Value = 1.723478864147651
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1961621590051053
$ws.Range("K23").Value = 0.491724344973818
$ws.Range("M23").Value = 0.2686257057968433
$ws.Range("O23").Value = 4.322100164057559

$ws.Range("B24").Value = 0.4737009348388597
$ws.Range("C24").Value = 0.04294444563628019
$ws.Range("D24").Value = 0.170257804102377
$ws.Range("E24").Value = 0.1614517151286776
$ws.Range("F24").Value = 1.723743363127326
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1949472521858837
$ws.Range("K24").Value = 0.4298745317985038
$ws.Range("M24").Value = 0.2487348262474782
$ws.Range("O24").Value = 4.340825147233375

$ws.Range("B25").Value = 0.4084097131298279
$ws.Range("C25").Value = 0.03602203066505183
$ws.Range("D25").Value = 0.1632458097342351
$ws.Range("E25").Value = 0.158537141183686
$ws.Range("F25").Value = 1.727140232460371
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1939663935487275
$ws.Range("K25").Value = 0.3632358490208105
$ws.Range("M25").Value = 0.2275850375963913
$ws.Range("O25").Value = 4.369058140513857
